$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1979.6
$ws.Range("I31").Value = 474.5
$ws.Range("K31").Value = 1423.5
$ws.Range("M31").Value = -1193.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 17267
$ws.Range("I70").Value = 20420.4
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 61261.2
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -60991.2
$ws.Range("N70").Value = -5040

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 17267
$ws.Range("I73").Value = 20420.4
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 61261.2
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -60325.2
$ws.Range("N73").Value = -6372

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 14706711
$ws.Range("I96").Value = 31250494
$ws.Range("J96").Value = 1126.5555
$ws.Range("K96").Value = 93751482
$ws.Range("L96").Value = 3379.6665
$ws.Range("M96").Value = -93750109
$ws.Range("N96").Value = -6125.666499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2219.0625
$ws.Range("I100").Value = 1782.2727
$ws.Range("K100").Value = 1782.2727
$ws.Range("M100").Value = -1241.2727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 853.2963
$ws.Range("J129").Value = 895.4375
$ws.Range("L129").Value = 2686.3125
$ws.Range("N129").Value = -12686.3125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 17660.129
$ws.Range("I137").Value = 1355.05
$ws.Range("J137").Value = 47305.727
$ws.Range("K137").Value = 4065.15
$ws.Range("L137").Value = 141917.181
$ws.Range("M137").Value = -1515.15
$ws.Range("N137").Value = -147017.181

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2623.4902
$ws.Range("J138").Value = 2802.4473
$ws.Range("L138").Value = 8407.341899999999
$ws.Range("N138").Value = -18687.3419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2331.2068
$ws.Range("I74").Value = 2522.25
$ws.Range("J74").Value = 1906.6666
$ws.Range("K74").Value = 2522.25
$ws.Range("L74").Value = 1906.6666
$ws.Range("M74").Value = -1648.25
$ws.Range("N74").Value = -3654.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2331.2068
$ws.Range("I77").Value = 2522.25
$ws.Range("J77").Value = 1906.6666
$ws.Range("K77").Value = 12611.25
$ws.Range("L77").Value = 9533.333000000001
$ws.Range("M77").Value = -8243.25
$ws.Range("N77").Value = -18269.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2762
$ws.Range("I97").Value = 3202.5
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 3202.5
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -2706.5
$ws.Range("N97").Value = -1992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5937.375
$ws.Range("I102").Value = 4583.1665
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 4583.1665
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -2961.1665
$ws.Range("N102").Value = -13244

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2201.35
$ws.Range("I122").Value = 2125.3225
$ws.Range("J122").Value = 2463.2222
$ws.Range("K122").Value = 6375.967500000001
$ws.Range("L122").Value = 7389.6666
$ws.Range("M122").Value = -3925.967500000001
$ws.Range("N122").Value = -12289.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 32407.295
$ws.Range("I132").Value = 2612.4443
$ws.Range("J132").Value = 65926.5
$ws.Range("K132").Value = 7837.3329
$ws.Range("L132").Value = 197779.5
$ws.Range("M132").Value = -5307.3329
$ws.Range("N132").Value = -202839.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2499.5
$ws.Range("I99").Value = 2499.5
$ws.Range("K99").Value = 2499.5
$ws.Range("M99").Value = -1001.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2001517.2
$ws.Range("I105").Value = 1371.5385
$ws.Range("K105").Value = 1371.5385
$ws.Range("M105").Value = 375.4614999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 48996.668
$ws.Range("J132").Value = 48996.668
$ws.Range("L132").Value = 48996.668
$ws.Range("N132").Value = -59116.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14963.788
$ws.Range("I31").Value = 27837.4
$ws.Range("J31").Value = 4235.778
$ws.Range("K31").Value = 27837.4
$ws.Range("L31").Value = 4235.778
$ws.Range("M31").Value = -27542.4
$ws.Range("N31").Value = -4825.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14963.788
$ws.Range("I34").Value = 27837.4
$ws.Range("J34").Value = 4235.778
$ws.Range("K34").Value = 27837.4
$ws.Range("L34").Value = 4235.778
$ws.Range("M34").Value = -27635.4
$ws.Range("N34").Value = -4639.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 533.3
$ws.Range("I107").Value = 547.875
$ws.Range("K107").Value = 547.875
$ws.Range("M107").Value = 1372.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 494.57895
$ws.Range("J5").Value = 525.05554
$ws.Range("L5").Value = 1575.16662
$ws.Range("N5").Value = -1799.16662

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6703.3335
$ws.Range("J68").Value = 9650.916999999999
$ws.Range("L68").Value = 28952.751
$ws.Range("N68").Value = -30574.751

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 6703.3335
$ws.Range("J71").Value = 9650.916999999999
$ws.Range("L71").Value = 86858.253
$ws.Range("N71").Value = -94970.253

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 3166.125
$ws.Range("J93").Value = 3554.8333
$ws.Range("L93").Value = 10664.4999
$ws.Range("N93").Value = -14408.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 139707.03
$ws.Range("I131").Value = 703.75
$ws.Range("J131").Value = 157082.44
$ws.Range("K131").Value = 2111.25
$ws.Range("L131").Value = 471247.32
$ws.Range("M131").Value = 2928.75
$ws.Range("N131").Value = -481327.32

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 494.57895
$ws.Range("J135").Value = 525.05554
$ws.Range("L135").Value = 4725.49986
$ws.Range("N135").Value = -9795.49986

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 61.92857
$ws.Range("J2").Value = 73.333336
$ws.Range("L2").Value = 73.333336
$ws.Range("N2").Value = -299.333336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4199.2856
$ws.Range("I97").Value = 1879
$ws.Range("K97").Value = 1879
$ws.Range("M97").Value = -1383

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2443.3
$ws.Range("I122").Value = 2290.5334
$ws.Range("K122").Value = 6871.600199999999
$ws.Range("M122").Value = -4421.600199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 49999
$ws.Range("J130").Value = 49999
$ws.Range("L130").Value = 49999
$ws.Range("N130").Value = -60039

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 79358.5
$ws.Range("I132").Value = 75870.92999999999
$ws.Range("J132").Value = 87496.164
$ws.Range("K132").Value = 227612.79
$ws.Range("L132").Value = 262488.492
$ws.Range("M132").Value = -225082.79
$ws.Range("N132").Value = -267548.492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2518.5386
$ws.Range("I93").Value = 2674.1
$ws.Range("K93").Value = 2674.1
$ws.Range("M93").Value = -1426.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5700
$ws.Range("I96").Value = 1750
$ws.Range("K96").Value = 1750
$ws.Range("M96").Value = -377
